$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 203, shifting existing rows 203:216 down to 204:217.
$ws.Range("A203").EntireRow.Insert()

# Populate the new row 203 with the weekly price entry (constant columns
# copied from the surrounding "Poroto granado" records).
$ws.Range("A203").Value = 5
$ws.Range("B203").Value = "Macroferia Regional de Talca"
$ws.Range("C203").Value = "Maule"
$ws.Range("D203").Value = 45021
$ws.Range("E203").Value = 7
$ws.Range("F203").Value = 100112030
$ws.Range("G203").Value = "Poroto granado"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 200
$ws.Range("K203").Value = 35000
$ws.Range("L203").Value = 35000
$ws.Range("M203").Value = 35000
$ws.Range("N203").Value = "`$/saco 25 kilos"
$ws.Range("O203").Value = "Región del Maule"
$ws.Range("P203").Value = 1400
$ws.Range("Q203").Value = 25
$ws.Range("R203").Value = "Hortaliza"
